$wb = $excel.ActiveWorkbook

# --- Add the new "missing_values" worksheet after the last sheet ("5d") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws5.Name = "missing_values"

# Header row
$ws5.Cells.Item(1,1).Value = "age"
$ws5.Cells.Item(1,2).Value = "sex\time"
$ws5.Cells.Item(1,3).Value = 2007
$ws5.Cells.Item(1,4).Value = 2010
$ws5.Cells.Item(1,5).Value = 2013

# Data rows (same content as sheet "3d" but with the age=1/H and age=4/F
# combinations removed -> "missing values")
$ws5.Cells.Item(2,1).Value = 0
$ws5.Cells.Item(2,2).Value = "F"
$ws5.Cells.Item(2,3).Value = 3722
$ws5.Cells.Item(2,4).Value = 3395
$ws5.Cells.Item(2,5).Value = 3347

$ws5.Cells.Item(3,1).Value = 0
$ws5.Cells.Item(3,2).Value = "H"
$ws5.Cells.Item(3,3).Value = 338
$ws5.Cells.Item(3,4).Value = 316
$ws5.Cells.Item(3,5).Value = 323

$ws5.Cells.Item(4,1).Value = 1
$ws5.Cells.Item(4,2).Value = "F"
$ws5.Cells.Item(4,3).Value = 2878
$ws5.Cells.Item(4,4).Value = 2791
$ws5.Cells.Item(4,5).Value = 2822

$ws5.Cells.Item(5,1).Value = 2
$ws5.Cells.Item(5,2).Value = "F"
$ws5.Cells.Item(5,3).Value = 4073
$ws5.Cells.Item(5,4).Value = 4161
$ws5.Cells.Item(5,5).Value = 4429

$ws5.Cells.Item(6,1).Value = 2
$ws5.Cells.Item(6,2).Value = "H"
$ws5.Cells.Item(6,3).Value = 1561
$ws5.Cells.Item(6,4).Value = 1463
$ws5.Cells.Item(6,5).Value = 1467

$ws5.Cells.Item(7,1).Value = 3
$ws5.Cells.Item(7,2).Value = "F"
$ws5.Cells.Item(7,3).Value = 3507
$ws5.Cells.Item(7,4).Value = 3741
$ws5.Cells.Item(7,5).Value = 3366

$ws5.Cells.Item(8,1).Value = 3
$ws5.Cells.Item(8,2).Value = "H"
$ws5.Cells.Item(8,3).Value = 2052
$ws5.Cells.Item(8,4).Value = 2052
$ws5.Cells.Item(8,5).Value = 2118

$ws5.Cells.Item(9,1).Value = 4
$ws5.Cells.Item(9,2).Value = "H"
$ws5.Cells.Item(9,3).Value = 3785
$ws5.Cells.Item(9,4).Value = 3508
$ws5.Cells.Item(9,5).Value = 3172

# --- Restore / set per-sheet selections (leftover cursor positions) ---
$ws2 = $wb.Worksheets.Item("2d")
$ws2.Range("C8").Select() | Out-Null

$ws3 = $wb.Worksheets.Item("3d")
$ws3.Range("B15").Select() | Out-Null

# Make the newly added sheet the active / selected tab, matching the
# author's final view state.
$ws5.Activate() | Out-Null
$ws5.Range("G24").Select() | Out-Null
